$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update Run No. (column D) and Bug (column F) for rows 4, 6, 8
$ws.Range("D4").Value = 3
$ws.Range("F4").Value = "Report not visible for FO"

$ws.Range("D6").Value = 3
$ws.Range("F6").Value = "Report not visible for FO"

$ws.Range("D8").Value = 3
$ws.Range("F8").Value = "Report not visible for FO"

# Update the selected range to reflect where the user left the cursor
$ws.Range("F8").Select()
